$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.018.84'
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").Value = '3.125.44'
$ws.Range("E3").Value = '  -5.62%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''521.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.57%  '
$ws.Range("D6").Value = '''134.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.08%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.121.06'
$ws.Range("E8").Value = '  -5.83%  '
$ws.Range("D9").Value = '''0.441'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.64%  '
$ws.Range("D10").Value = '''7.26'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.76%  '
$ws.Range("E11").Value = '  -8.77%  '
$ws.Range("E12").Value = '  -6.84%  '
$ws.Range("D13").Value = '3.652.66'
$ws.Range("E13").Value = '  -5.94%  '
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").Value = '''25.45'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.47%  '
$ws.Range("D16").Value = '3.119.90'
$ws.Range("E16").Value = '  -5.81%  '
$ws.Range("D17").Value = '57.906.38'
$ws.Range("E17").Value = '  -4.16%  '
$ws.Range("E18").Value = '  -8.26%  '
$ws.Range("D19").Value = '''5.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.60%  '
$ws.Range("E20").Value = '  -9.54%  '
$ws.Range("D21").Value = '''7.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.40%  '
$ws.Range("D22").Value = '''343.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.28%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '''68.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.67%  '
$ws.Range("D25").Value = '''0.507'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.95%  '
$ws.Range("D26").Value = '3.251.84'
$ws.Range("E26").Value = '  -5.59%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0₃0955'
$ws.Range("E27").Value = '  -6.72%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.166'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.61%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '''6.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.73%  '
$ws.Range("E32").Value = '  -9.34%  '
$ws.Range("D33").Value = '''21.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.56%  '
$ws.Range("D34").Value = '''6.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.56%  '
$ws.Range("E35").Value = '  -2.08%  '
$ws.Range("D36").Value = '''157.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.91%  '
$ws.Range("D37").Value = '''4.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.78%  '
$ws.Range("D38").Value = '''6.17'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.18%  '
$ws.Range("D39").Value = '''1.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.15%  '
$ws.Range("E40").Value = '  -5.84%  '
$ws.Range("B41").Value = 'RenzoRestakedETH'
$ws.Range("C41").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D41").Value = '3.148.55'
$ws.Range("E41").Value = '  -5.81%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''40.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''23.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -11.11%  '
$ws.Range("D44").Value = '''0.694'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.75%  '
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("D46").Value = '''3.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.68%  '
$ws.Range("D47").Value = '''0.998'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '''1.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.46%  '
$ws.Range("D49").Value = '2.255.64'
$ws.Range("E49").Value = '  -5.01%  '
$ws.Range("D50").Value = '''6.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.57%  '
$ws.Range("D51").Value = '''20.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.66%  '
